$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of row 3 with row 4, and row 5 with row 6
# (record identity/order changed upstream), while keeping row 2 and
# everything below row 6 untouched.

# ---- New row 3 (previously row 4's data) ----
$ws.Range("A3").Value = 130861152
$ws.Range("B3").Value = 91805
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 1108
$ws.Range("F3").Value = "Harticka"
$ws.Range("G3").Value = "Pelloporus leporinus"
$ws.Range("H3").Value = "(Fr.) Krieglst."
$ws.Range("P3").Value = "Djupbäcken, Jmt"
$ws.Range("Q3").Value = 442868
$ws.Range("R3").Value = 7039767
$ws.Range("S3").Value = 10
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("AC3").Value = "I stående levande gran med full längd."
$ws.Range("AW3").Value = "Kristian Zackrisson"
$ws.Range("AX3").Value = "Kristian Zackrisson"

# ---- New row 4 (previously row 3's data) ----
$ws.Range("A4").Value = 130853761
$ws.Range("B4").Value = 79245
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 230405
$ws.Range("F4").Value = "Garnlav (ssp. sarmentosa)"
$ws.Range("G4").Value = "Alectoria sarmentosa subsp. sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("P4").Value = "Djupbäcken, Djupbäcken, Jmt"
$ws.Range("Q4").Value = 442771
$ws.Range("R4").Value = 7039709
$ws.Range("S4").Value = 20
$ws.Range("Z4").Value = "11:05"
$ws.Range("AB4").Value = "11:05"
$ws.Range("AC4").ClearContents()
$ws.Range("AW4").Value = "Maria Danvind"
$ws.Range("AX4").Value = "Maria Danvind"

# ---- New row 5 (previously row 6's data) ----
$ws.Range("A5").Value = 130861156
$ws.Range("Q5").Value = 442897
$ws.Range("R5").Value = 7039676
$ws.Range("AC5").Value = "På död stående gran med full längd."

# ---- New row 6 (previously row 5's data) ----
$ws.Range("A6").Value = 130861158
$ws.Range("Q6").Value = 442743
$ws.Range("R6").Value = 7039650
$ws.Range("AC6").Value = "På gran."
